$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title.
# ------------------------------------------------------------------
$metaPara = Find-ParagraphByText $d "Meta description"
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Free Aztec Warrior Slot Game |
#    DragonGaming Review" right before the final (italic "Prompt: ...")
#    paragraph -- i.e. right after "Not ideal for high rollers".
# ------------------------------------------------------------------
$anchorPara = Find-ParagraphByText $d "Not ideal for high rollers"
$rBreak = $anchorPara.Range
$rBreak.Collapse(0)
$rBreak.InsertAfter("`r")

$promptPara = Find-ParagraphByText $d "Prompt: Create a feature image"
$rNew = $d.Range($anchorPara.Range.End, $promptPara.Range.Start)
$rNewContent = $d.Range($rNew.Start, $rNew.End - 1)
$rNewContent.Text = "Play Free Aztec Warrior Slot Game | DragonGaming Review"
$rNewContent.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Replace the old "Prompt: ..." text (last paragraph) with the meta
#    description sentence, keeping the run's italic formatting intact.
# ------------------------------------------------------------------
$promptPara = Find-ParagraphByText $d "Prompt: Create a feature image"
$rLast = $promptPara.Range
$rLastContent = $d.Range($rLast.Start, $rLast.End - 1)
$rLastContent.Text = "Experience the Aztec civilization with DragonGaming's Aztec Warrior free slot game. Read our review and play for free today!"
